# Fruta / hortaliza, semanal
# Insert one new week of Kiwi price records (Mercado Mayorista Lo Valledor de
# Santiago) at the top of the data block that starts at row 1045, pushing the
# existing rows (old 1045-1078) down by four rows (new 1049-1082).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows right before the current row 1045. Excel copies the
# formatting (incl. the date number format on column D) from the row above,
# exactly like the existing rows in this block.
$ws.Range("A1045:A1048").EntireRow.Insert()

# Columns that are constant for every record in this sheet/block.
$marketId   = 6
$market     = "Mercado Mayorista Lo Valledor de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$tipo       = "Fruta"
$productId  = 100101
$producto   = "Berries"
$categoriaId = 100101007
$categoria  = "Kiwi"
$variedad   = "Hayward"
$unidad     = "`$/bins (450 kilos)"
$origen     = "Región de O'Higgins"
$kgUnidad   = 450
$fecha      = 45041

# New rows: [row, calidad, volumen, precioMin, precioMax, precioProm, precioKg]
$newRows = @(
    @(1045, "Especial",               25, 300000, 300000, 300000, 667),
    @(1046, "Extra (doble especial)",  18, 350000, 350000, 350000, 778),
    @(1047, "Primera",                 39, 250000, 270000, 259231, 576),
    @(1048, "Segunda",                 17, 220000, 220000, 220000, 489)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value2  = $marketId
    $ws.Cells.Item($r, 2).Value   = $market
    $ws.Cells.Item($r, 3).Value   = $region
    $ws.Cells.Item($r, 4).Value2  = $fecha
    $ws.Cells.Item($r, 5).Value2  = $codreg
    $ws.Cells.Item($r, 6).Value   = $tipo
    $ws.Cells.Item($r, 7).Value2  = $productId
    $ws.Cells.Item($r, 8).Value   = $producto
    $ws.Cells.Item($r, 9).Value2  = $categoriaId
    $ws.Cells.Item($r, 10).Value  = $categoria
    $ws.Cells.Item($r, 11).Value  = $variedad
    $ws.Cells.Item($r, 12).Value  = $row[1]
    $ws.Cells.Item($r, 13).Value2 = $row[2]
    $ws.Cells.Item($r, 14).Value2 = $row[3]
    $ws.Cells.Item($r, 15).Value2 = $row[4]
    $ws.Cells.Item($r, 16).Value2 = $row[5]
    $ws.Cells.Item($r, 17).Value  = $unidad
    $ws.Cells.Item($r, 18).Value  = $origen
    $ws.Cells.Item($r, 19).Value2 = $row[6]
    $ws.Cells.Item($r, 20).Value2 = $kgUnidad
}
